{"js": "// Word JavaScript API (Office.js) edit script.\n//\n// Two independent changes, mirroring the target diff:\n//  1. In the \"l|ine , cast , through, bounce, utmost, setback\" paragraph,\n//     drop the (Word-managed) \"_GoBack\" bookmark right after \"through\" -\n//     this simply tracks where the author's cursor last was after\n//     editing, and it moves away from wherever it previously sat.\n//  2. In the trailing numbered vocabulary list (items \"1.\" - \"9.\"), the\n//     answers were shuffled around: \"intrusion\" (1) and \"tabloid\" (2)\n//     moved down to (5) and (8) respectively, while \"coverage\" moved up\n//     from (8) to (1), and \"allegation\" moved up from (5) to (2). Items\n//     3, 4, 6, 7 and 9 stay exactly where they were. Because \"_GoBack\" is\n//     unique, it gets removed from its old spot (right after \"scandal\",\n//     at the very end of the document) once it is re-added next to\n//     \"through\".\n\n// ---------------------------------------------------------------------\n// Change 1 (part a): drop the old \"_GoBack\" bookmark (it currently sits\n// right after \"scandal\", at the very end of the document) before adding\n// the new one, since Word only ever keeps a single \"_GoBack\" bookmark.\n// ---------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 1 (part b): split \"...through, bounce...\" and drop the\n// \"_GoBack\" bookmark right after \"through\".\n// ---------------------------------------------------------------------\nconst headlineParagraph = paragraphs.items[2]; // \"line , cast , through, bounce, utmost, setback\"\nconst throughResults = headlineParagraph.search(\"through\", { matchCase: true });\nthroughResults.load(\"items\");\nawait context.sync();\n\nconst throughRange = throughResults.items[0];\nconst afterThrough = throughRange.getRange(\"End\");\nafterThrough.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Change 2: reorder the answers in the numbered list (paragraphs 38-45,\n// i.e. items \"1.\" through \"8.\" - \"9.\" does not change).\n// ---------------------------------------------------------------------\nconst wordSwaps = [\n  { paragraphIndex: 38, oldWord: \"intrusion\", newWord: \"coverage\" },   // 1. intrusion  -> coverage\n  { paragraphIndex: 39, oldWord: \"tabloid\", newWord: \"allegation\" },   // 2. tabloid    -> allegation\n  { paragraphIndex: 42, oldWord: \"allegation\", newWord: \"intrusion\" }, // 5. allegation -> intrusion\n  { paragraphIndex: 45, oldWord: \"coverage\", newWord: \"tabloid\" },     // 8. coverage   -> tabloid\n];\n\nfor (const swap of wordSwaps) {\n  const freshParagraphs = body.paragraphs;\n  freshParagraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = freshParagraphs.items[swap.paragraphIndex];\n  const found = paragraph.search(swap.oldWord, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n\n  found.items[0].insertText(swap.newWord, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Two independent changes, mirroring the target diff:\n#  1. In the \"l|ine , cast , through, bounce, utmost, setback\" paragraph,\n#     drop the (Word-managed) \"_GoBack\" bookmark right after \"through\" -\n#     this simply tracks where the author's cursor last was after\n#     editing, and it moves away from wherever it previously sat.\n#  2. In the trailing numbered vocabulary list (items \"1.\" - \"9.\"), the\n#     answers were shuffled around: \"intrusion\" (1) and \"tabloid\" (2)\n#     moved down to (5) and (8) respectively, while \"coverage\" moved up\n#     from (8) to (1), and \"allegation\" moved up from (5) to (2). Items\n#     3, 4, 6, 7 and 9 stay exactly where they were. Because \"_GoBack\" is\n#     unique, it gets removed from its old spot (right after \"scandal\",\n#     at the very end of the document) once it is re-added next to\n#     \"through\".\n\n$d = $word.ActiveDocument\n\n# -----------------------------------------------------------------------\n# Change 1 (part a): drop the old \"_GoBack\" bookmark (it currently sits\n# right after \"scandal\", at the very end of the document) before adding\n# the new one, since Word only ever keeps a single \"_GoBack\" bookmark.\n# -----------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# -----------------------------------------------------------------------\n# Change 1 (part b): split \"...through, bounce...\" and drop the\n# \"_GoBack\" bookmark right after \"through\".\n# -----------------------------------------------------------------------\n$headlineParagraph = $d.Paragraphs.Item(3).Range  # \"line , cast , through, bounce, utmost, setback\"\n$headlineText = $headlineParagraph.Text\n$throughOffset = $headlineText.IndexOf(\"through\")\n$throughEnd = $headlineParagraph.Start + $throughOffset + \"through\".Length\n$bookmarkRange = $d.Range($throughEnd, $throughEnd)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n\n# -----------------------------------------------------------------------\n# Change 2: reorder the answers in the numbered list (paragraphs 39-46\n# in this 1-indexed collection, i.e. items \"1.\" through \"8.\" - \"9.\" does\n# not change).\n# -----------------------------------------------------------------------\nfunction Replace-WordInParagraph($paragraphIndex, $oldWord, $newWord) {\n    $paragraphRange = $d.Paragraphs.Item($paragraphIndex).Range\n    $paragraphText = $paragraphRange.Text\n    $wordOffset = $paragraphText.IndexOf($oldWord)\n    $wordStart = $paragraphRange.Start + $wordOffset\n    $wordEnd = $wordStart + $oldWord.Length\n    $targetRange = $d.Range($wordStart, $wordEnd)\n    $targetRange.Text = $newWord\n}\n\nReplace-WordInParagraph 39 \"intrusion\" \"coverage\"     # 1. intrusion  -> coverage\nReplace-WordInParagraph 40 \"tabloid\" \"allegation\"      # 2. tabloid    -> allegation\nReplace-WordInParagraph 43 \"allegation\" \"intrusion\"    # 5. allegation -> intrusion\nReplace-WordInParagraph 46 \"coverage\" \"tabloid\"        # 8. coverage   -> tabloid\n"}
